$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column C (shifts Mutual Fund..QoQ one column right)
$ws.Range("C1").EntireColumn.Insert()

# Set the header value for the newly inserted column
$ws.Range("C1").Value = "Industry"

# Copy the style/formatting from the adjacent header cell (B1) so the new
# header matches the look (bold, centered, bordered) of the other headers
$ws.Range("B1").Copy()
$ws.Range("C1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
